$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4354
$ws.Range("I6").Value = 4354
$ws.Range("K6").Value = 13062
$ws.Range("M6").Value = -12950
$ws.Range("H132").Value = 219581.16
$ws.Range("I132").Value = 238858.84
$ws.Range("J132").Value = 79129.42999999999
$ws.Range("K132").Value = 716576.52
$ws.Range("L132").Value = 237388.29
$ws.Range("M132").Value = -714046.52
$ws.Range("N132").Value = -242448.29

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 41205.6
$ws.Range("I2").Value = 51108.85
$ws.Range("K2").Value = 51108.85
$ws.Range("M2").Value = -50995.85
$ws.Range("H116").Value = 41205.6
$ws.Range("I116").Value = 51108.85
$ws.Range("K116").Value = 51108.85
$ws.Range("M116").Value = -48814.85
$ws.Range("H122").Value = 2156.2693
$ws.Range("I122").Value = 1871.9231
$ws.Range("J122").Value = 2440.6155
$ws.Range("K122").Value = 5615.7693
$ws.Range("L122").Value = 7321.8465
$ws.Range("M122").Value = -3165.7693
$ws.Range("N122").Value = -12221.8465
$ws.Range("H132").Value = 1930.6389
$ws.Range("I132").Value = 1402.1724
$ws.Range("K132").Value = 4206.5172
$ws.Range("M132").Value = -1676.5172

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 41205.6
$ws.Range("I3").Value = 51108.85
$ws.Range("K3").Value = 51108.85
$ws.Range("M3").Value = -50994.85

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1145.4706
$ws.Range("I5").Value = 459.6316
$ws.Range("J5").Value = 2014.2
$ws.Range("K5").Value = 1378.8948
$ws.Range("L5").Value = 6042.6
$ws.Range("M5").Value = -1266.8948
$ws.Range("N5").Value = -6266.6
$ws.Range("H7").Value = 180.2
$ws.Range("I7").Value = 101
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 303
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -191
$ws.Range("N7").Value = -824
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H116").Value = 1371.2858
$ws.Range("I116").Value = 1099.8334
$ws.Range("K116").Value = 3299.5002
$ws.Range("M116").Value = 142.4998000000001
$ws.Range("H122").Value = 562.9583
$ws.Range("I122").Value = 261.86667
$ws.Range("J122").Value = 1064.7778
$ws.Range("K122").Value = 2356.80003
$ws.Range("L122").Value = 9583.0002
$ws.Range("M122").Value = 93.19997000000012
$ws.Range("N122").Value = -14483.0002
$ws.Range("H135").Value = 1145.4706
$ws.Range("I135").Value = 459.6316
$ws.Range("J135").Value = 2014.2
$ws.Range("K135").Value = 4136.6844
$ws.Range("L135").Value = 18127.8
$ws.Range("M135").Value = -1601.6844
$ws.Range("N135").Value = -23197.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 3321.6316
$ws.Range("I126").Value = 4250
$ws.Range("J126").Value = 3074.0667
$ws.Range("K126").Value = 12750
$ws.Range("L126").Value = 9222.2001
$ws.Range("M126").Value = -10280
$ws.Range("N126").Value = -14162.2001
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 44999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 44999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 44999
$ws.Range("N129").Value = -54999
$ws.Range("H130").Value = 49800
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 49800
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 49800
$ws.Range("N130").Value = -59840
$ws.Range("H131").Value = 35000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 35000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080
$ws.Range("H132").Value = 5126.6
$ws.Range("I132").Value = 5518.6
$ws.Range("J132").Value = 4342.6
$ws.Range("K132").Value = 16555.8
$ws.Range("L132").Value = 13027.8
$ws.Range("M132").Value = -14025.8
$ws.Range("N132").Value = -18087.8
$ws.Range("H133").Value = 19997.143
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 19997.143
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 19997.143
$ws.Range("N133").Value = -30117.143
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 142891710
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 142891710
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 142891710
$ws.Range("N135").Value = -142901850
$ws.Range("H136").Value = 21774.8
$ws.Range("I136").Value = 22222
$ws.Range("J136").Value = 21663
$ws.Range("K136").Value = 66666
$ws.Range("L136").Value = 64989
$ws.Range("M136").Value = -64116
$ws.Range("N136").Value = -70089
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 64500
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 64500
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 64500
$ws.Range("N138").Value = -74780
$ws.Range("H139").Value = 35000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 78429
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 78429
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 78429
$ws.Range("N141").Value = -88789

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2880.32
$ws.Range("I7").Value = 2088.5
$ws.Range("J7").Value = 3252.9412
$ws.Range("K7").Value = 2088.5
$ws.Range("L7").Value = 3252.9412
$ws.Range("M7").Value = -1976.5
$ws.Range("N7").Value = -3476.9412
$ws.Range("H122").Value = 3719.1904
$ws.Range("I122").Value = 3433.8333
$ws.Range("K122").Value = 10301.4999
$ws.Range("M122").Value = -7851.499899999999
$ws.Range("H126").Value = 2880.32
$ws.Range("I126").Value = 2088.5
$ws.Range("J126").Value = 3252.9412
$ws.Range("K126").Value = 6265.5
$ws.Range("L126").Value = 9758.8236
$ws.Range("M126").Value = -3795.5
$ws.Range("N126").Value = -14698.8236

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 835693.7
$ws.Range("I81").Value = 1251696.8
$ws.Range("J81").Value = 3687.75
$ws.Range("K81").Value = 2503393.6
$ws.Range("L81").Value = 7375.5
$ws.Range("M81").Value = -2502332.6
$ws.Range("N81").Value = -9497.5
$ws.Range("H84").Value = 835693.7
$ws.Range("I84").Value = 1251696.8
$ws.Range("J84").Value = 3687.75
$ws.Range("K84").Value = 12516968
$ws.Range("L84").Value = 36877.5
$ws.Range("M84").Value = -12511664
$ws.Range("N84").Value = -47485.5
$ws.Range("H128").Value = 42552
$ws.Range("J128").Value = 42552
$ws.Range("L128").Value = 42552
$ws.Range("N128").Value = -52512
